$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  "Name",
  "Odell Beckham Jr. - NYG",
  "Antonio Brown - PIT",
  "DeAndre Hopkins - HOU",
  "Julio Jones - ATL",
  "Keenan Allen - LAC",
  "Michael Thomas - NO",
  "Davante Adams - GB",
  "T.Y. Hilton - IND",
  "A.J. Green - CIN",
  "Stefon Diggs - MIN",
  "Mike Evans - TB",
  "Larry Fitzgerald - ARI",
  "Doug Baldwin - SEA",
  "Adam Thielen - MIN",
  "Amari Cooper - OAK",
  "Emmanuel Sanders - DEN",
  "Brandin Cooks - LAR",
  "Sammy Watkins - KC",
  "Josh Gordon - CLE",
  "Allen Robinson - CHI",
  "Demaryius Thomas - DEN",
  "Alshon Jeffery - PHI",
  "Kelvin Benjamin - BUF",
  "Golden Tate - DET",
  "Marvin Jones - DET",
  "Randall Cobb - GB",
  "Jarvis Landry - CLE",
  "JuJu Smith-Schuster - PIT",
  "Corey Davis - TEN",
  "Chris Hogan - NE",
  "Robby Anderson - NYJ",
  "Jamison Crowder - WAS",
  "Cameron Meredith - NO",
  "DeSean Jackson - TB",
  "Mike Williams - LAC",
  "Cooper Kupp - LAR",
  "Will Fuller - HOU",
  "Marquise Goodwin - SF",
  "Pierre Garcon - SF",
  "Allen Hurns - DAL",
  "Robert Woods - LAR",
  "Taylor Gabriel - CHI",
  "Dede Westbrook - JAC",
  "D.J. Moore - CAR",
  "Devin Funchess - CAR",
  "Michael Crabtree - BAL",
  "Josh Doctson - WAS",
  "Dez Bryant - ",
  "Jordy Nelson - OAK",
  "Martavis Bryant - OAK",
  "Calvin Ridley - ATL",
  "Sterling Shepard - NYG",
  "Julian Edelman - NE",
  "Michael Gallup - DAL",
  "Christian Kirk - ARI",
  "Devante Parker - MIA",
  "Kenny Stills - MIA",
  "Rishard Matthews - TEN",
  "Nelson Agholor - PHI",
  "Marqise Lee - JAC",
  "Kenny Golladay - DET",
  "Tyler Lockett - SEA",
  "Tyrell Williams - LAC",
  "Mohamed Sanu - ATL",
  "Paul Richardson - WAS",
  "Ted Ginn - NO",
  "Anthony Miller - CHI",
  "Chris Godwin - TB",
  "Jordan Matthews - ",
  "Donte Moncrief - JAC",
  "Keelan Cole - JAC",
  "Quincy Enunwa - NYJ",
  "Corey Coleman - BUF",
  "Jermaine Kearse - NYJ",
  "Danny Amendola - MIA",
  "John Brown - BAL",
  "Geronimo Allison - GB",
  "Terrance Williams - DAL",
  "Mike Wallace - PHI",
  "Brandon Marshall - SEA",
  "Albert Wilson - MIA",
  "Ryan Grant - IND",
  "Cole Beasley - DAL",
  "Terrelle Pryor - NYJ",
  "Dante Pettis - SF",
  "Courtland Sutton - DEN",
  "Chester Rogers - IND",
  "James Washington - PIT",
  "Willie Snead - BAL",
  "Travis Benjamin - LAC",
  "Brice Butler - ARI",
  "J.J. Nelson - ARI",
  "Taywan Taylor - TEN",
  "Zay Jones - BUF",
  "John Ross - CIN",
  "Trent Taylor - SF",
  "Brandon LaFell - ",
  "Torrey Smith - CAR",
  "Kendall Wright - MIN",
  "Jeremy Maclin - ",
  "Tavon Austin - DAL",
  "Curtis Samuel - CAR",
  "Jaron Brown - SEA",
  "Adam Humphries - TB",
  "Chad Williams - ARI",
  "Kevin White - CHI",
  "Eric Decker - NE",
  "Tre'Quan Smith - NO",
  "Malcolm Mitchell - ",
  "Jakeem Grant - MIA",
  "Keke Coutee - HOU",
  "J'Mon Moore - GB",
  "Kenny Britt - NE",
  "Equanimeous St. Brown - GB",
  "Chris Conley - KC",
  "Bruce Ellington - HOU",
  "Amara Darboh - SEA",
  "Laquon Treadwell - MIN",
  "D.J. Chark - JAC",
  "Tyler Boyd - CIN",
  "Jeremy Kerley - BUF",
  "Carlos Henderson - DEN",
  "Mack Hollins - PHI",
  "Dontrelle Inman - ",
  "Josh Reynolds - LAR",
  "Deon Cain - IND",
  "Cordarrelle Patterson - NE",
  "Brandon Coleman - ",
  "Demarcus Robinson - KC",
  "Seth Roberts - OAK",
  "Roger Lewis - NYG",
  "Cody Latimer - NYG",
  "Andre Holmes - BUF",
  "Antonio Callaway - CLE",
  "Braxton Miller - HOU",
  "DaeSean Hamilton - DEN",
  "Chris Moore - BAL",
  "TJ Jones - DET",
  "Aldrick Robinson - SF",
  "Breshad Perriman - BAL",
  "Chad Hansen - NYJ",
  "Justin Hardy - ATL",
  "Tajae Sharpe - TEN",
  "Daurice Fountain - IND",
  "Kamar Aiken - PHI",
  "Jarius Wright - CAR",
  "Ryan Switzer - OAK",
  "Josh Malone - CIN",
  "ArDarius Stewart - NYJ",
  "Marcell Ateman - OAK",
  "Jordan Lasley - BAL",
  "Jaleel Scott - BAL",
  "Deonte Thompson - DAL",
  "Kendrick Bourne - SF",
  "Johnny Holton - OAK",
  "Darrius Heyward-Bey - PIT",
  "Rod Streater - BUF",
  "Geremy Davis - LAC",
  "Pharoh Cooper - LAR",
  "Michael Campanaro - TEN",
  "Rashard Higgins - CLE",
  "Maurice Harris - WAS",
  "De'Anthony Thomas - KC",
  "Phillip Dorsett - NE",
  "Josh Bellamy - CHI",
  "Justin Watson - TB",
  "Tanner McEvoy - SEA",
  "Brian Quick - WAS",
  "Marquez Valdes-Scantling - GB",
  "Ricardo Louis - CLE",
  "Justin Hunter - PIT",
  "Leonte Carroo - MIA",
  "Jaydon Mickens - JAC",
  "Cedrick Wilson - DAL",
  "Auden Tate - CIN",
  "Trevor Davis - GB",
  "Simmie Cobbs Jr. - WAS",
  "Deontay Burnett - TEN",
  "Marvin Hall - ATL",
  "Andy Jones - DET",
  "Alex Erickson - CIN",
  "Travis Rudolph - NYG",
  "Jehu Chesson - KC",
  "Nick Williams - TEN",
  "Bennie Fowler - CHI",
  "Marcus Johnson - SEA",
  "Ray-Ray McCloud - BUF",
  "Trey Quinn - WAS",
  "Damiere Byrd - CAR",
  "Tavarres King - MIN",
  "Jordan Taylor - DEN",
  "Russell Shepard - NYG",
  "Bobo Wilson - TB",
  "Tommylee Lewis - NO",
  "Chris Thompson - HOU",
  "Artavis Scott - LAC",
  "Kaelin Clay - BUF",
  "Markus Wheaton - PHI",
  "Mike Thomas - LAR",
  "Stacy Coley - MIN",
  "Dylan Cantrell - LAC",
  "Bradley Marquez - DET",
  "Sammie Coates - HOU",
  "Kasen Williams - IND",
  "Freddie Martino - TB",
  "Damion Ratley - CLE",
  "Dwayne Harris - OAK",
  "Michael Floyd - NO",
  "Shelton Gibson - PHI",
  "Braxton Berrios - NE",
  "Noah Brown - DAL",
  "DeAngelo Yancey - GB",
  "Krishawn Hogan - IND",
  "Jace Billingsley - DET",
  "Cody Core - CIN",
  "Andre Roberts - NYJ",
  "Robert Davis - WAS",
  "Bryce Treggs - PHI",
  "Russell Gage - ATL",
  "Marquess Wilson - ",
  "Michael Clark - ",
  "Jeff Janis - CLE",
  "Tanner Gentry - CHI",
  "Richie James - SF",
  "Devin Fuller - ",
  "Griff Whalen - OAK",
  "Isaiah McKenzie - DEN",
  "Daniel Braverman - KC",
  "Rashad Greene - JAC",
  "Cobi Hamilton - ",
  "Jake Wieneke - MIN",
  "Jaelen Strong - ",
  "Reggie Davis - ATL",
  "Austin Carr - NO",
  "Charone Peake - NYJ",
  "Quincy Adeboyejo - BAL",
  "DeVier Posey - BAL",
  "Austin Proehl - BUF",
  "Riley McCarron - NE",
  "David Moore - SEA",
  "K.J. Brent - IND",
  "KD Cannon - DAL",
  "Aaron Burbridge - SF",
  "Zach Pascal - IND",
  "Devin Smith - ",
  "Brandon Reilly - BUF",
  "Hunter Sharp - NYG"
)

for ($i = 0; $i -lt $values.Count; $i++) {
  $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

$ws.Cells.Item(249, 1).Value = ""

$ws.Columns.Item(1).ColumnWidth = 26.6
$ws.Range("A1").Select()
